# "Add today's walk in Athens"
# The walk-tracking workbook keeps a running "today's walk" distance in
# Sheet1!G2 (added onto yesterday's running total in F1 to produce the new
# running total in F2 via the existing formula =F1+G2). Record today's walk
# by updating G2 - F2 recalculates automatically from its existing formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("G2").Value = 54.3

# Recalculate so the dependent formula (F2 = F1+G2) is fresh before saving.
$excel.Calculate()

# Leave the selection where the user's cursor landed after entering the value.
$ws.Range("H3").Select() | Out-Null
